$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P1").Value = "擬調利率"

$ws.Range("P5").Select()
